$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Insert a new "Vorbedingung" bullet after
#    "Es müssen bereits Nährstoffe im System vorhanden sein."
#    announcing that the recipe dialog is already open (TeamC note).
# ------------------------------------------------------------------
$anchorIdx = Find-ParagraphIndex("*Nährstoffe im System vorhanden sein*")
if ($anchorIdx -gt 0) {
    $anchorPara = $d.Paragraphs($anchorIdx)
    $insertRange = $anchorPara.Range
    $insertRange.Collapse(0)            # wdCollapseEnd
    $insertRange.InsertParagraphAfter()

    # The freshly-created paragraph inherits the preceding paragraph's
    # numbering / run formatting (Arial list item), so we only need to
    # set its text.
    $newPara = $d.Paragraphs($anchorIdx + 1)
    $newPara.Range.Text = "Dialog für das erstellen ein Rezept ist bereits geöffnet.(TeamC)"
}

# ------------------------------------------------------------------
# 2) Remove the "Zubereitung (Pflichtfeld)" bullet (and its bookmark)
#    entirely.
# ------------------------------------------------------------------
$targetIdx = Find-ParagraphIndex("Zubereitung (Pflichtfeld)*")
if ($targetIdx -gt 0) {
    $targetPara = $d.Paragraphs($targetIdx)
    $targetPara.Range.Delete()
}
